$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "98.842.74"
$ws.Range("E2").Value = "  +1.46%  "
$ws.Range("D3").Value = "3.397.47"
$ws.Range("E3").Value = "  +8.53%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "'261.07"
$ws.Range("E5").Value = "  +8.55%  "
$ws.Range("D6").Value = "'637.59"
$ws.Range("E6").Value = "  +4.63%  "
$ws.Range("E7").Value = "  +25.78%  "
$ws.Range("D8").Value = "'0.398"
$ws.Range("E8").Value = "  +2.72%  "
$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("D10").Value = "'0.892"
$ws.Range("E10").Value = "  +12.02%  "
$ws.Range("D11").Value = "3.397.69"
$ws.Range("E11").Value = "  +8.66%  "
$ws.Range("D12").Value = "'0.201"
$ws.Range("E12").Value = "  +1.63%  "
$ws.Range("D13").Value = "98.585.26"
$ws.Range("E13").Value = "  +2.04%  "
$ws.Range("D14").Value = "'36.44"
$ws.Range("E14").Value = "  +6.96%  "
$ws.Range("E15").Value = "  +3.93%  "
$ws.Range("D16").Value = "3.987.14"
$ws.Range("E16").Value = "  +7.64%  "
$ws.Range("D17").Value = "'5.59"
$ws.Range("E17").Value = "  +4.19%  "
$ws.Range("D18").Value = "3.401.34"
$ws.Range("E18").Value = "  +9.31%  "
$ws.Range("D19").Value = "'3.64"
$ws.Range("E19").Value = "  +1.31%  "
$ws.Range("D20").Value = "'15.30"
$ws.Range("E20").Value = "  +5.15%  "
$ws.Range("D21").Value = "'495.79"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").Value = "'6.23"
$ws.Range("E22").Value = "  +9.38%  "
$ws.Range("E23").Value = "  +8.65%  "
$ws.Range("D24").Value = "'9.44"
$ws.Range("E24").Value = "  +6.87%  "
$ws.Range("D25").Value = "'5.76"
$ws.Range("E25").Value = "  +4.02%  "
$ws.Range("D26").Value = "'89.31"
$ws.Range("E26").Value = "  +3.56%  "
$ws.Range("D27").Value = "'12.16"
$ws.Range("E27").Value = "  +4.18%  "
$ws.Range("E29").Value = "  +20.47%  "
$ws.Range("D30").Value = "'0.998"
$ws.Range("E30").Value = "  -0.14%  "
$ws.Range("D31").Value = "'0.196"
$ws.Range("E31").Value = "  +11.16%  "
$ws.Range("D32").Value = "'0.132"
$ws.Range("E32").Value = "  +5.51%  "
$ws.Range("D33").Value = "'9.67"
$ws.Range("E33").Value = "  +6.76%  "
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  +17.84%  "
$ws.Range("D35").Value = "'28.06"
$ws.Range("E35").Value = "  +6.51%  "
$ws.Range("D36").Value = "'7.42"
$ws.Range("E36").Value = "  +0.48%  "
$ws.Range("E37").Value = "  +6.72%  "
$ws.Range("E38").Value = "  -0.34%  "
$ws.Range("D39").Value = "'0.472"
$ws.Range("E39").Value = "  +7.14%  "
$ws.Range("D40").Value = "'507.79"
$ws.Range("E40").Value = "  +3.30%  "
$ws.Range("D41").Value = "'24.87"
$ws.Range("E41").Value = "  +2.75%  "
$ws.Range("D42").Value = "'3.75"
$ws.Range("E42").Value = "  +2.92%  "
$ws.Range("D43").Value = "'1.28"
$ws.Range("E43").Value = "  +3.24%  "
$ws.Range("D45").Value = "'0.789"
$ws.Range("E45").Value = "  +12.09%  "
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").Value = "'160.41"
$ws.Range("E47").Value = "  -0.57%  "
$ws.Range("D48").Value = "'1.95"
$ws.Range("E48").Value = "  +1.64%  "
$ws.Range("D49").Value = "'4.71"
$ws.Range("E49").Value = "  +7.88%  "
$ws.Range("D50").Value = "'46.75"
$ws.Range("E50").Value = "  +5.06%  "
$ws.Range("E51").Value = "  +13.12%  "
